# Weekly price-sheet update: a new observation (week of 2023-09-08) is
# inserted as row 27, pushing all the existing data rows (old rows 27-107)
# down by one (new rows 28-108). The sheet's used range grows from
# A1:T107 to A1:T108.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 27, shifting rows 27-107
# down to 28-108 (and the sheet dimension grows to A1:T108 accordingly).
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly observation.
$ws.Cells.Item(27, 1).Value  = 3
$ws.Cells.Item(27, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(27, 3).Value  = "Coquimbo"
$ws.Cells.Item(27, 4).Value  = 45177
$ws.Cells.Item(27, 5).Value  = 5
$ws.Cells.Item(27, 6).Value  = "Fruta"
$ws.Cells.Item(27, 7).Value  = 100108
$ws.Cells.Item(27, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(27, 9).Value  = 100108004
$ws.Cells.Item(27, 10).Value = "Papaya"
$ws.Cells.Item(27, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(27, 12).Value = "Segunda"
$ws.Cells.Item(27, 13).Value = 50
$ws.Cells.Item(27, 14).Value = 17000
$ws.Cells.Item(27, 15).Value = 17000
$ws.Cells.Item(27, 16).Value = 17000
$ws.Cells.Item(27, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(27, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(27, 19).Value = 1700
$ws.Cells.Item(27, 20).Value = 10

# Make sure the new date cell carries the same date number-format style
# as the rest of column D (Excel's native row-insert already copies this
# from the row above, but set it explicitly to be safe).
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
